# Update the Fitness column (C) values for run_25 log rows 2-183
# (Generation 0 through 181), matching the recorded GA run output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12998
$ws.Range("C3:C6").Value = 11207
$ws.Range("C7:C8").Value = 10671
$ws.Range("C9:C11").Value = 10533
$ws.Range("C12:C18").Value = 8591
$ws.Range("C19:C23").Value = 8167
$ws.Range("C24:C53").Value = 7828
$ws.Range("C54:C72").Value = 7594
$ws.Range("C73:C183").Value = 7569
